$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Editorial corrections to affiliations / surnames
$ws.Range("C13").Value = "Leibniz Institute for Prevention Research and Epidemiology - BIPS, Bremen"
$ws.Range("C23").Value = "Swiss Tropical and Public Health Institute"
$ws.Range("B27").Value = "Sorenson"

# Update the active selection to reflect where the editor ended up
$ws.Range("B27").Select()
